$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.011.11'
$ws.Range('E2').Value = '  +6.18%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.593.52'
$ws.Range('E3').Value = '  +6.28%  '
$ws.Range('E4').Value = '  +0.43%  '
$ws.Range('E5').Value = '  +3.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.46'
$ws.Range('E6').Value = '  +8.62%  '
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.547'
$ws.Range('E8').Value = '  +3.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.621.86'
$ws.Range('E9').Value = '  +7.39%  '
$ws.Range('E10').Value = '  +6.44%  '
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('E12').Value = '  +3.01%  '
$ws.Range('E13').Value = '  +4.87%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.67'
$ws.Range('E14').Value = '  +4.11%  '
$ws.Range('E15').Value = '  +6.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.063.15'
$ws.Range('E16').Value = '  +6.17%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '65.801.13'
$ws.Range('E17').Value = '  +5.65%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.628.80'
$ws.Range('E18').Value = '  +7.75%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.22'
$ws.Range('E19').Value = '  +6.74%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.24'
$ws.Range('E20').Value = '  +5.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '354.33'
$ws.Range('E21').Value = '  +10.86%  '
$ws.Range('E22').Value = '  +5.28%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.25'
$ws.Range('E23').Value = '  +5.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.998'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.07'
$ws.Range('E25').Value = '  +2.65%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '66.26'
$ws.Range('E26').Value = '  +2.27%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '642.15'
$ws.Range('E27').Value = '  +0.80%  '
$ws.Range('E28').Value = '  +13.52%  '
$ws.Range('E29').Value = '  +6.19%  '
$ws.Range('E30').Value = '  +8.96%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.990'
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('E32').Value = '  +6.10%  '
$ws.Range('E33').Value = '  +5.79%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.138'
$ws.Range('E34').Value = '  +6.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.64'
$ws.Range('E35').Value = '  +11.69%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.995'
$ws.Range('E36').Value = '  -0.33%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.01'
$ws.Range('E37').Value = '  +8.75%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.96'
$ws.Range('E38').Value = '  +9.99%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.68'
$ws.Range('E39').Value = '  +9.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.44'
$ws.Range('E40').Value = '  +6.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '156.10'
$ws.Range('E41').Value = '  +3.77%  '
$ws.Range('E42').Value = '  +3.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.83'
$ws.Range('E43').Value = '  +8.65%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.42'
$ws.Range('E44').Value = '  +1.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0₆0317'
$ws.Range('E45').Value = '  +4.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '163.28'
$ws.Range('E46').Value = '  +7.93%  '
$ws.Range('E47').Value = '  -0.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.20'
$ws.Range('E48').Value = '  +6.00%  '
$ws.Range('E49').Value = '  +7.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '21.93'
$ws.Range('E50').Value = '  +10.07%  '
$ws.Range('E51').Value = '  +6.72%  '
